$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (string) number format on the Price and Volume(1h) columns that are
# being updated, so Excel does not reinterpret the values as numbers/percentages.
$priceVolumeCells = @(
    "D2",
    "E2",
    "D3",
    "E3",
    "D4",
    "E4",
    "D5",
    "E5",
    "D6",
    "E6",
    "D7",
    "E7",
    "D8",
    "E8",
    "D9",
    "E9",
    "D10",
    "E10",
    "D11",
    "E11",
    "D12",
    "E12",
    "E13",
    "D14",
    "E14",
    "D15",
    "E15",
    "D16",
    "E16",
    "D17",
    "E17",
    "D18",
    "E19",
    "D20",
    "E20",
    "D21",
    "E21",
    "D22",
    "E22",
    "D23",
    "E23",
    "D24",
    "E24",
    "D25",
    "E25",
    "D26",
    "E26",
    "D27",
    "E27",
    "D39",
    "E39",
    "D40",
    "E40",
    "D41",
    "E41",
    "D42",
    "E42",
    "D43",
    "E43",
    "D44",
    "E44",
    "D45",
    "E45",
    "D46",
    "E46",
    "D47",
    "E47"
)
foreach ($addr in $priceVolumeCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values cell by cell, row by row, matching the refreshed
# coinranking.com crypto symbol list snapshot.

# Row 2
$ws.Range("D2").Value = "296.46"
$ws.Range("E2").Value = "-4.94%"

# Row 3
$ws.Range("D3").Value = "31.50"
$ws.Range("E3").Value = "-2.79%"

# Row 4
$ws.Range("D4").Value = "5.115"
$ws.Range("E4").Value = "-4.11%"

# Row 5
$ws.Range("D5").Value = "0.07467"
$ws.Range("E5").Value = "-2.39%"

# Row 6
$ws.Range("D6").Value = "7.735"
$ws.Range("E6").Value = "-1.68%"

# Row 7
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "3.793"
$ws.Range("E7").Value = "2.19%"

# Row 8
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").Value = "1.661"
$ws.Range("E8").Value = "3.62%"

# Row 9
$ws.Range("D9").Value = "0.9323"
$ws.Range("E9").Value = "1.32%"

# Row 10
$ws.Range("D10").Value = "0.1694"
$ws.Range("E10").Value = "-2.15%"

# Row 11
$ws.Range("D11").Value = "0.07146"
$ws.Range("E11").Value = "-5.92%"

# Row 12
$ws.Range("D12").Value = "0.07979"
$ws.Range("E12").Value = "-3.09%"

# Row 13
$ws.Range("E13").Value = "0.38%"

# Row 14
$ws.Range("D14").Value = "0.09913"
$ws.Range("E14").Value = "0.33%"

# Row 15
$ws.Range("D15").Value = "0.001495"
$ws.Range("E15").Value = "-1.60%"

# Row 16
$ws.Range("D16").Value = "0.006272"
$ws.Range("E16").Value = "2.60%"

# Row 17
$ws.Range("D17").Value = "3.461"
$ws.Range("E17").Value = "-0.78%"

# Row 18
$ws.Range("D18").Value = "2.224"

# Row 19
$ws.Range("E19").Value = "-0.96%"

# Row 20
$ws.Range("D20").Value = "0.1336"
$ws.Range("E20").Value = "1.34%"

# Row 21
$ws.Range("D21").Value = "4.570"
$ws.Range("E21").Value = "7.00%"

# Row 22
$ws.Range("D22").Value = "0.04653"
$ws.Range("E22").Value = "2.15%"

# Row 23
$ws.Range("D23").Value = "0.1553"
$ws.Range("E23").Value = "-4.57%"

# Row 24
$ws.Range("D24").Value = "0.001218"
$ws.Range("E24").Value = "-0.18%"

# Row 25
$ws.Range("D25").Value = "0.004422"
$ws.Range("E25").Value = "-1.69%"

# Row 26
$ws.Range("D26").Value = "0.0001303"
$ws.Range("E26").Value = "0.42%"

# Row 27
$ws.Range("D27").Value = "0.0001878"
$ws.Range("E27").Value = "7.97%"

# Row 39
$ws.Range("D39").Value = "0.01672"
$ws.Range("E39").Value = "-1.61%"

# Row 40
$ws.Range("D40").Value = "0.04448"
$ws.Range("E40").Value = "-3.79%"

# Row 41
$ws.Range("D41").Value = "0.007070"
$ws.Range("E41").Value = "-2.50%"

# Row 42
$ws.Range("D42").Value = "0.1328"
$ws.Range("E42").Value = "-3.09%"

# Row 43
$ws.Range("D43").Value = "0.002105"
$ws.Range("E43").Value = "-6.68%"

# Row 44
$ws.Range("D44").Value = "0.01224"
$ws.Range("E44").Value = "-14.26%"

# Row 45
$ws.Range("D45").Value = "0.00005990"
$ws.Range("E45").Value = "-3.44%"

# Row 46
$ws.Range("D46").Value = "0.7116"
$ws.Range("E46").Value = "-62.40%"

# Row 47
$ws.Range("D47").Value = "0.01101"
$ws.Range("E47").Value = "-15.20%"
